$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.65656533333333
$ws.Range("H2").Value = 67.969696
$ws.Range("I2").Value = 0.9268638682343595
$ws.Range("J2").Value = 0.9268638682343595
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.445876
$ws.Range("N2").Value = 1.337628
$ws.Range("O2").Value = 0.004821482820471209
$ws.Range("P2").Value = 0.00482148282047121
$ws.Range("Q2").Value = 10.10201872456533
$ws.Range("R2").Value = 90.91816852108801
$ws.Range("S2").Value = 0.004468858217607455
$ws.Range("T2").Value = 0.004468858217607456

$ws.Range("G3").Value = 22.65656533333333
$ws.Range("H3").Value = 67.969696
$ws.Range("I3").Value = 0.9268638682343595
$ws.Range("J3").Value = 0.9268638682343595
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 85.56602466666668
$ws.Range("N3").Value = 256.698074
$ws.Range("O3").Value = 0.925268724816651
$ws.Range("P3").Value = 0.9252687248166511
$ws.Range("Q3").Value = 1938.632228173945
$ws.Range("R3").Value = 17447.6900535655
$ws.Range("S3").Value = 0.8575981494398343
$ws.Range("T3").Value = 0.8575981494398344

$ws.Range("G4").Value = 22.65656533333333
$ws.Range("H4").Value = 67.969696
$ws.Range("I4").Value = 0.9268638682343595
$ws.Range("J4").Value = 0.9268638682343595
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.465044
$ws.Range("N4").Value = 19.395132
$ws.Range("O4").Value = 0.06990979236287771
$ws.Range("P4").Value = 0.06990979236287773
$ws.Range("Q4").Value = 146.4756917688747
$ws.Range("R4").Value = 1318.281225919872
$ws.Range("S4").Value = 0.06479686057691772
$ws.Range("T4").Value = 0.06479686057691773

$ws.Range("G5").Value = 0.9818753333333333
$ws.Range("H5").Value = 2.945626
$ws.Range("I5").Value = 0.04016781697437198
$ws.Range("J5").Value = 0.04016781697437198
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.445876
$ws.Range("N5").Value = 1.337628
$ws.Range("O5").Value = 0.004821482820471209
$ws.Range("P5").Value = 0.00482148282047121
$ws.Range("Q5").Value = 0.4377946461253333
$ws.Range("R5").Value = 3.940151815128
$ws.Range("S5").Value = 0.0001936684394777663
$ws.Range("T5").Value = 0.0001936684394777664

$ws.Range("G6").Value = 0.9818753333333333
$ws.Range("H6").Value = 2.945626
$ws.Range("I6").Value = 0.04016781697437198
$ws.Range("J6").Value = 0.04016781697437198
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 85.56602466666668
$ws.Range("N6").Value = 256.698074
$ws.Range("O6").Value = 0.925268724816651
$ws.Range("P6").Value = 0.9252687248166511
$ws.Range("Q6").Value = 84.01516899159157
$ws.Range("R6").Value = 756.136520924324
$ws.Range("S6").Value = 0.03716602479054579
$ws.Range("T6").Value = 0.03716602479054579

$ws.Range("G7").Value = 0.9818753333333333
$ws.Range("H7").Value = 2.945626
$ws.Range("I7").Value = 0.04016781697437198
$ws.Range("J7").Value = 0.04016781697437198
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.465044
$ws.Range("N7").Value = 19.395132
$ws.Range("O7").Value = 0.06990979236287771
$ws.Range("P7").Value = 0.06990979236287773
$ws.Range("Q7").Value = 6.347867232514666
$ws.Range("R7").Value = 57.130805092632
$ws.Range("S7").Value = 0.00280812374434842
$ws.Range("T7").Value = 0.002808123744348421

$ws.Range("G8").Value = 0.8058883333333333
$ws.Range("H8").Value = 2.417665
$ws.Range("I8").Value = 0.03296831479126849
$ws.Range("J8").Value = 0.03296831479126849
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.445876
$ws.Range("N8").Value = 1.337628
$ws.Range("O8").Value = 0.004821482820471209
$ws.Range("P8").Value = 0.00482148282047121
$ws.Range("Q8").Value = 0.3593262665133333
$ws.Range("R8").Value = 3.23393639862
$ws.Range("S8").Value = 0.0001589561633859879
$ws.Range("T8").Value = 0.0001589561633859879

$ws.Range("G9").Value = 0.8058883333333333
$ws.Range("H9").Value = 2.417665
$ws.Range("I9").Value = 0.03296831479126849
$ws.Range("J9").Value = 0.03296831479126849
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 85.56602466666668
$ws.Range("N9").Value = 256.698074
$ws.Range("O9").Value = 0.925268724816651
$ws.Range("P9").Value = 0.9252687248166511
$ws.Range("Q9").Value = 68.95666100857889
$ws.Range("R9").Value = 620.60994907721
$ws.Range("S9").Value = 0.03050455058627093
$ws.Range("T9").Value = 0.03050455058627093

$ws.Range("G10").Value = 0.8058883333333333
$ws.Range("H10").Value = 2.417665
$ws.Range("I10").Value = 0.03296831479126849
$ws.Range("J10").Value = 0.03296831479126849
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.465044
$ws.Range("N10").Value = 19.395132
$ws.Range("O10").Value = 0.06990979236287771
$ws.Range("P10").Value = 0.06990979236287773
$ws.Range("Q10").Value = 5.210103534086667
$ws.Range("R10").Value = 46.89093180678
$ws.Range("S10").Value = 0.00230480804161157
$ws.Range("T10").Value = 0.00230480804161157

